$d = $word.ActiveDocument

$pairs = @(
    @("871×5=", "925×9="),
    @("682×9=", "287×4="),
    @("732×5=", "992×5="),
    @("951×2=", "212×7="),
    @("806×5=", "930×3="),
    @("795×5=", "846×9="),
    @("755×5=", "260×3="),
    @("222×6=", "165×3="),
    @("577×4=", "325×4="),
    @("910×5=", "821×4="),
    @("805×2=", "413×2="),
    @("434×8=", "935×3="),
    @("744×7=", "133×6="),
    @("130×6=", "696×3="),
    @("852×8=", "621×2="),
    @("880×3=", "296×6="),
    @("710×4=", "766×4="),
    @("947×3=", "832×5="),
    @("387×3=", "782×6="),
    @("639×6=", "940×7="),
    @("754×7=", "312×8="),
    @("284×8=", "984×9="),
    @("826×3=", "776×2="),
    @("932×7=", "818×9="),
    @("596×6=", "798×4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
